{"js": "// Apply the \"Added many more features\" edit to the Gold Factory review.\n// Each change is a straightforward text replacement; we locate the old\n// text with body.search() (exact, case-sensitive) and replace the whole\n// hit with the new text, which preserves the run's existing formatting.\n\nconst replacements = [\n  // Title heading (Heading1) AND the bold \"call to action\" line near the\n  // end share identical original text, so search() will find both and we\n  // replace every hit.\n  {\n    find: \"Play Gold Factory Slot Game for Free\",\n    replace: \"Play Gold Factory Slot For Free\",\n  },\n  {\n    find: \"User-friendly layout and engaging graphics\",\n    replace: \"User-friendly and engaging layout\",\n  },\n  {\n    find: \"Possible maximum payout of 619,000 coins\",\n    replace: \"Well-designed graphics and sound effects\",\n  },\n  {\n    find: \"Well-placed sound effects\",\n    replace: \"Attractive jackpot and bonus games\",\n  },\n  {\n    find: \"Linear gameplay may not be for everyone\",\n    replace: \"Limited betting options with a maximum bet of \\u20AC1,000\",\n  },\n  {\n    find: \"RTP of 95.6% is average compared to other slots\",\n    replace: \"RTP of 95.6% is relatively average\",\n  },\n  {\n    find: \"Explore the imaginary factory of the Gold Factory slot game and win big by playing for free. Enjoy exciting bonus features and a user-friendly layout.\",\n    replace: \"Read our review of Gold Factory slot game and try it out for free.\",\n  },\n];\n\nfor (const { find, replace } of replacements) {\n  const results = context.document.body.search(find, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${find}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(replace, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the \"Added many more features\" edit to the Gold Factory review.\n# Each change is a straightforward text replacement performed with\n# Find/Replace over the whole document (wdReplaceAll), which preserves\n# the existing run formatting of the matched text.\n\n$d = $word.ActiveDocument\n\nfunction Replace-All($findText, $replaceText) {\n    $r = $d.Content\n    $r.Find.ClearFormatting()\n    $found = $r.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n    if (-not $found) {\n        throw \"Text not found: $findText\"\n    }\n}\n\n# Title heading (Heading1) AND the bold \"call to action\" line near the end\n# share identical original text, so ReplaceAll updates both occurrences.\nReplace-All \"Play Gold Factory Slot Game for Free\" \"Play Gold Factory Slot For Free\"\n\nReplace-All \"User-friendly layout and engaging graphics\" \"User-friendly and engaging layout\"\n\nReplace-All \"Possible maximum payout of 619,000 coins\" \"Well-designed graphics and sound effects\"\n\nReplace-All \"Well-placed sound effects\" \"Attractive jackpot and bonus games\"\n\nReplace-All \"Linear gameplay may not be for everyone\" \"Limited betting options with a maximum bet of \u20ac1,000\"\n\nReplace-All \"RTP of 95.6% is average compared to other slots\" \"RTP of 95.6% is relatively average\"\n\nReplace-All \"Explore the imaginary factory of the Gold Factory slot game and win big by playing for free. Enjoy exciting bonus features and a user-friendly layout.\" \"Read our review of Gold Factory slot game and try it out for free.\"\n"}
